# Add season-record columns (Wins / Losses / Ties) to the DET_2019 sheet.
# Every row gets the team's season record: 47 wins, 114 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used by the rest of the header row (bold,
# bordered, centered) by copying the format from the neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-55): season record for every player row ---
$firstDataRow = 2
$lastDataRow = 55
$rowCount = $lastDataRow - $firstDataRow + 1

$records = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $records[$i,0] = 47
    $records[$i,1] = 114
    $records[$i,2] = 0
}

$ws.Range("AD" + $firstDataRow + ":AF" + $lastDataRow).Value = $records
